$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A slightly (OOXML stored width goes from 15.42578125 to 15.7109375)
$ws.Columns.Item(1).ColumnWidth = 14.8

# Replace the 33 sampled values in column A with the new dataset

$ws.Range("A1").Value = 0.34985991912974157
$ws.Range("A2").Value = -0.092103155724920782
$ws.Range("A3").Value = -0.0089999997585223923
$ws.Range("A4").Value = -0.011999999936993078
$ws.Range("A5").Value = -0.0059999997639712532
$ws.Range("A6").Value = -0.0059999997590907128
$ws.Range("A7").Value = -0.019999999712835503
$ws.Range("A8").Value = 0.00039299887388111188
$ws.Range("A9").Value = -0.0059999997572131036
$ws.Range("A10").Value = -0.0059999997566180241
$ws.Range("A11").Value = -0.0044999997614993958
$ws.Range("A12").Value = -0.005999999756632679
$ws.Range("A13").Value = -0.0059999997570825414
$ws.Range("A14").Value = -0.011999999737756006
$ws.Range("A15").Value = 0.043602941198907175
$ws.Range("A16").Value = -0.0059999997573352282
$ws.Range("A17").Value = -0.00599999975626897
$ws.Range("A18").Value = -0.0089999997461829295
$ws.Range("A19").Value = -0.0089999997641170282
$ws.Range("A20").Value = -0.0089999997618992467
$ws.Range("A21").Value = -0.0089999997615421989
$ws.Range("A22").Value = -0.0089999997612908444
$ws.Range("A23").Value = -0.008999999756190924
$ws.Range("A24").Value = -0.041999999643870467
$ws.Range("A25").Value = -0.041999999641814334
$ws.Range("A26").Value = -0.0059999997585933329
$ws.Range("A27").Value = -0.0059999997579582853
$ws.Range("A28").Value = -0.0059999997553568107
$ws.Range("A29").Value = -0.011999999734022992
$ws.Range("A30").Value = -0.019999999707202232
$ws.Range("A31").Value = 0.025120672619896922
$ws.Range("A32").Value = -0.02099999970339983
$ws.Range("A33").Value = -0.0059999997524515791
